$d = $word.ActiveDocument

# 1. Update the italic subtitle under "Missing Persons Outlier Detection"
$rng = $d.Content
$null = $rng.Find.Execute("Statistical Anomaly Detection for Trafficking & Organized Crime")
$rng.Text = "Geospatial Crime Pattern Analysis | 41,200 NamUs Cases"

# 2. Update the first bullet (41,200 cases -> 7 statistical methods + 3 ML models)
$rng = $d.Content
$null = $rng.Find.Execute("Analyzed 41,200 cases across 101 years identifying trafficking corridors at up to 46.86σ significance")
$rng.Text = "7 statistical methods + 3 ML models detecting trafficking corridors at up to 46.86σ significance"

# 3. Insert a new bullet after the one we just updated (I-35 corridor acceleration)
$rng = $d.Content
$null = $rng.Find.Execute("7 statistical methods + 3 ML models detecting trafficking corridors at up to 46.86σ significance")
$para = $rng.Paragraphs(1)
$null = $para.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($para.Index + 1)
$insertPoint = $d.Range($newPara.Range.Start, $newPara.Range.Start)
$insertPoint.InsertBefore("• I-35 corridor acceleration: 170% increase in missing persons, structural break at 2020")

# 4. Update the Streamlit dashboard bullet
$rng = $d.Content
$null = $rng.Find.Execute("Built 7-page interactive Streamlit dashboard with geospatial visualization")
$rng.Text = "Live Streamlit dashboard with spatial autocorrelation (Moran's I), ARIMA forecasting, and LISA clustering"

Write-Host "Edits applied"
